# Generate Report for Handoff
# Updates the localization-status report: flips the "In Translation" status
# to "Ready for handoff" and refreshes the related generation timestamps,
# then re-widens the Status columns to fit the new (longer) text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet --------------------------------------------------
# zh-cn / de-de status columns
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
# Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-09-07 06:57:28"

# --- zh-cn sheet -------------------------------------------------------
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-07 06:57:22"

# --- de-de sheet -------------------------------------------------------
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-07 06:57:28"

# --- Re-fit the Status columns to the new, longer text ----------------
# (ColumnWidth is expressed in characters; the sheet's stored XML width
# is ColumnWidth + 0.8333... so we back out the character width that
# yields the target display width.)
$newStatusWidth = 17.2159881591797 - 0.8333333333333334

$wsOverview.Columns.Item(5).ColumnWidth = $newStatusWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusWidth

Write-Host "Report regenerated for handoff."
